$d = $word.ActiveDocument

# Locate the paragraph that contains "Baz chan" (it is split across two runs
# around a _GoBack bookmark: "Baz chan" + bookmark + "ges").
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Baz chan*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $r = $target.Range

    # Replace the leading "Baz chan" text with the new sentence.
    $r.Find.Execute("Baz chan", $true, $false, $false, $false, $false, $true, 1, $false,
                     "Version management is important but I hope I never have to use GitHub again.", 2)

    # Remove the trailing "ges" run that followed the bookmark, scoped to this
    # paragraph only so it can't touch "changes" elsewhere in the document.
    $r2 = $target.Range
    $r2.Find.Execute("ges", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
}
